$wb = $excel.ActiveWorkbook

# Rename the sheet "HistoriadelapoblacióndeAsi" to "Historia"
$ws = $wb.Worksheets.Item("HistoriadelapoblacióndeAsi")
$ws.Name = "Historia"

# The chart on "GraficaHistoria" references this sheet by name in its
# series formulas; update them to use the new sheet name.
$gws = $wb.Worksheets.Item("GraficaHistoria")
$co = $gws.ChartObjects(1)
$chart = $co.Chart
$chart.SeriesCollection(1).Formula = "=SERIES(Historia!`$A`$76,Historia!`$A`$2:`$A`$76,Historia!`$A`$2:`$A`$76,1)"
$chart.SeriesCollection(2).Formula = "=SERIES(Historia!`$B`$1,,Historia!`$B`$2:`$B`$76,2)"

# Update the selected cell on the "Historia" sheet from F22 to E20
$ws.Activate()
$ws.Range("E20").Select()
